$d = $word.ActiveDocument

$d.Content.Find.Execute("33×69=", $true, $false, $false, $false, $false, $true, 1, $false, "50×51=", 2) | Out-Null
$d.Content.Find.Execute("22×61=", $true, $false, $false, $false, $false, $true, 1, $false, "23×47=", 2) | Out-Null
$d.Content.Find.Execute("96×39=", $true, $false, $false, $false, $false, $true, 1, $false, "68×91=", 2) | Out-Null
$d.Content.Find.Execute("85×20=", $true, $false, $false, $false, $false, $true, 1, $false, "100×93=", 2) | Out-Null
$d.Content.Find.Execute("63×16=", $true, $false, $false, $false, $false, $true, 1, $false, "69×69=", 2) | Out-Null
$d.Content.Find.Execute("86×22=", $true, $false, $false, $false, $false, $true, 1, $false, "99×21=", 2) | Out-Null
$d.Content.Find.Execute("31×32=", $true, $false, $false, $false, $false, $true, 1, $false, "76×26=", 2) | Out-Null
$d.Content.Find.Execute("80×64=", $true, $false, $false, $false, $false, $true, 1, $false, "13×65=", 2) | Out-Null
$d.Content.Find.Execute("27×95=", $true, $false, $false, $false, $false, $true, 1, $false, "52×56=", 2) | Out-Null
$d.Content.Find.Execute("76×22=", $true, $false, $false, $false, $false, $true, 1, $false, "63×22=", 2) | Out-Null
$d.Content.Find.Execute("58×75=", $true, $false, $false, $false, $false, $true, 1, $false, "92×84=", 2) | Out-Null
$d.Content.Find.Execute("77×84=", $true, $false, $false, $false, $false, $true, 1, $false, "41×50=", 2) | Out-Null
$d.Content.Find.Execute("28×23=", $true, $false, $false, $false, $false, $true, 1, $false, "72×93=", 2) | Out-Null
$d.Content.Find.Execute("85×72=", $true, $false, $false, $false, $false, $true, 1, $false, "39×31=", 2) | Out-Null
$d.Content.Find.Execute("53×60=", $true, $false, $false, $false, $false, $true, 1, $false, "55×71=", 2) | Out-Null
$d.Content.Find.Execute("74×58=", $true, $false, $false, $false, $false, $true, 1, $false, "92×25=", 2) | Out-Null
$d.Content.Find.Execute("40×94=", $true, $false, $false, $false, $false, $true, 1, $false, "49×33=", 2) | Out-Null
$d.Content.Find.Execute("20×36=", $true, $false, $false, $false, $false, $true, 1, $false, "11×48=", 2) | Out-Null
$d.Content.Find.Execute("74×13=", $true, $false, $false, $false, $false, $true, 1, $false, "100×39=", 2) | Out-Null
$d.Content.Find.Execute("70×99=", $true, $false, $false, $false, $false, $true, 1, $false, "44×79=", 2) | Out-Null
$d.Content.Find.Execute("48×30=", $true, $false, $false, $false, $false, $true, 1, $false, "47×64=", 2) | Out-Null
$d.Content.Find.Execute("17×45=", $true, $false, $false, $false, $false, $true, 1, $false, "44×49=", 2) | Out-Null
$d.Content.Find.Execute("74×82=", $true, $false, $false, $false, $false, $true, 1, $false, "66×49=", 2) | Out-Null
$d.Content.Find.Execute("70×11=", $true, $false, $false, $false, $false, $true, 1, $false, "74×66=", 2) | Out-Null
$d.Content.Find.Execute("72×59=", $true, $false, $false, $false, $false, $true, 1, $false, "99×69=", 2) | Out-Null
$d.Content.Find.Execute("11×33=", $true, $false, $false, $false, $false, $true, 1, $false, "68×73=", 2) | Out-Null
$d.Content.Find.Execute("52×58=", $true, $false, $false, $false, $false, $true, 1, $false, "21×59=", 2) | Out-Null
$d.Content.Find.Execute("76×57=", $true, $false, $false, $false, $false, $true, 1, $false, "91×40=", 2) | Out-Null
$d.Content.Find.Execute("82×76=", $true, $false, $false, $false, $false, $true, 1, $false, "60×65=", 2) | Out-Null
$d.Content.Find.Execute("39×62=", $true, $false, $false, $false, $false, $true, 1, $false, "11×34=", 2) | Out-Null
$d.Content.Find.Execute("72×61=", $true, $false, $false, $false, $false, $true, 1, $false, "51×61=", 2) | Out-Null
$d.Content.Find.Execute("71×30=", $true, $false, $false, $false, $false, $true, 1, $false, "36×73=", 2) | Out-Null
$d.Content.Find.Execute("34×65=", $true, $false, $false, $false, $false, $true, 1, $false, "27×43=", 2) | Out-Null
$d.Content.Find.Execute("13×77=", $true, $false, $false, $false, $false, $true, 1, $false, "31×86=", 2) | Out-Null
$d.Content.Find.Execute("49×80=", $true, $false, $false, $false, $false, $true, 1, $false, "41×58=", 2) | Out-Null
$d.Content.Find.Execute("37×54=", $true, $false, $false, $false, $false, $true, 1, $false, "85×38=", 2) | Out-Null
$d.Content.Find.Execute("50×26=", $true, $false, $false, $false, $false, $true, 1, $false, "49×53=", 2) | Out-Null
$d.Content.Find.Execute("45×61=", $true, $false, $false, $false, $false, $true, 1, $false, "72×81=", 2) | Out-Null
$d.Content.Find.Execute("96×54=", $true, $false, $false, $false, $false, $true, 1, $false, "34×84=", 2) | Out-Null
$d.Content.Find.Execute("68×86=", $true, $false, $false, $false, $false, $true, 1, $false, "23×46=", 2) | Out-Null
$d.Content.Find.Execute("13×18=", $true, $false, $false, $false, $false, $true, 1, $false, "81×40=", 2) | Out-Null
$d.Content.Find.Execute("24×83=", $true, $false, $false, $false, $false, $true, 1, $false, "49×50=", 2) | Out-Null
$d.Content.Find.Execute("71×81=", $true, $false, $false, $false, $false, $true, 1, $false, "33×13=", 2) | Out-Null
$d.Content.Find.Execute("75×69=", $true, $false, $false, $false, $false, $true, 1, $false, "70×10=", 2) | Out-Null
$d.Content.Find.Execute("55×69=", $true, $false, $false, $false, $false, $true, 1, $false, "10×52=", 2) | Out-Null
$d.Content.Find.Execute("46×22=", $true, $false, $false, $false, $false, $true, 1, $false, "55×18=", 2) | Out-Null
$d.Content.Find.Execute("29×16=", $true, $false, $false, $false, $false, $true, 1, $false, "91×37=", 2) | Out-Null
$d.Content.Find.Execute("23×60=", $true, $false, $false, $false, $false, $true, 1, $false, "57×60=", 2) | Out-Null
$d.Content.Find.Execute("42×11=", $true, $false, $false, $false, $false, $true, 1, $false, "58×62=", 2) | Out-Null
$d.Content.Find.Execute("21×14=", $true, $false, $false, $false, $false, $true, 1, $false, "43×99=", 2) | Out-Null
$d.Content.Find.Execute("40×73=", $true, $false, $false, $false, $false, $true, 1, $false, "68×66=", 2) | Out-Null
$d.Content.Find.Execute("41×24=", $true, $false, $false, $false, $false, $true, 1, $false, "98×61=", 2) | Out-Null
$d.Content.Find.Execute("93×94=", $true, $false, $false, $false, $false, $true, 1, $false, "27×68=", 2) | Out-Null
$d.Content.Find.Execute("94×50=", $true, $false, $false, $false, $false, $true, 1, $false, "29×74=", 2) | Out-Null
$d.Content.Find.Execute("42×30=", $true, $false, $false, $false, $false, $true, 1, $false, "42×31=", 2) | Out-Null
$d.Content.Find.Execute("45×51=", $true, $false, $false, $false, $false, $true, 1, $false, "76×43=", 2) | Out-Null
$d.Content.Find.Execute("53×19=", $true, $false, $false, $false, $false, $true, 1, $false, "27×85=", 2) | Out-Null
$d.Content.Find.Execute("86×67=", $true, $false, $false, $false, $false, $true, 1, $false, "97×62=", 2) | Out-Null
$d.Content.Find.Execute("68×50=", $true, $false, $false, $false, $false, $true, 1, $false, "51×84=", 2) | Out-Null
$d.Content.Find.Execute("88×43=", $true, $false, $false, $false, $false, $true, 1, $false, "81×34=", 2) | Out-Null
$d.Content.Find.Execute("41×28=", $true, $false, $false, $false, $false, $true, 1, $false, "88×33=", 2) | Out-Null
$d.Content.Find.Execute("18×13=", $true, $false, $false, $false, $false, $true, 1, $false, "11×19=", 2) | Out-Null
$d.Content.Find.Execute("67×93=", $true, $false, $false, $false, $false, $true, 1, $false, "69×63=", 2) | Out-Null
$d.Content.Find.Execute("75×93=", $true, $false, $false, $false, $false, $true, 1, $false, "35×75=", 2) | Out-Null
$d.Content.Find.Execute("13×28=", $true, $false, $false, $false, $false, $true, 1, $false, "50×46=", 2) | Out-Null
$d.Content.Find.Execute("26×93=", $true, $false, $false, $false, $false, $true, 1, $false, "46×20=", 2) | Out-Null
$d.Content.Find.Execute("42×73=", $true, $false, $false, $false, $false, $true, 1, $false, "91×100=", 2) | Out-Null
$d.Content.Find.Execute("11×30=", $true, $false, $false, $false, $false, $true, 1, $false, "16×41=", 2) | Out-Null
$d.Content.Find.Execute("28×34=", $true, $false, $false, $false, $false, $true, 1, $false, "33×25=", 2) | Out-Null
$d.Content.Find.Execute("98×83=", $true, $false, $false, $false, $false, $true, 1, $false, "99×92=", 2) | Out-Null
$d.Content.Find.Execute("30×68=", $true, $false, $false, $false, $false, $true, 1, $false, "49×45=", 2) | Out-Null
$d.Content.Find.Execute("20×93=", $true, $false, $false, $false, $false, $true, 1, $false, "58×49=", 2) | Out-Null
$d.Content.Find.Execute("59×77=", $true, $false, $false, $false, $false, $true, 1, $false, "34×88=", 2) | Out-Null
$d.Content.Find.Execute("44×46=", $true, $false, $false, $false, $false, $true, 1, $false, "41×91=", 2) | Out-Null
$d.Content.Find.Execute("32×11=", $true, $false, $false, $false, $false, $true, 1, $false, "40×19=", 2) | Out-Null
$d.Content.Find.Execute("34×32=", $true, $false, $false, $false, $false, $true, 1, $false, "98×100=", 2) | Out-Null
$d.Content.Find.Execute("83×50=", $true, $false, $false, $false, $false, $true, 1, $false, "29×53=", 2) | Out-Null
$d.Content.Find.Execute("44×92=", $true, $false, $false, $false, $false, $true, 1, $false, "97×65=", 2) | Out-Null
$d.Content.Find.Execute("22×74=", $true, $false, $false, $false, $false, $true, 1, $false, "67×90=", 2) | Out-Null
$d.Content.Find.Execute("16×71=", $true, $false, $false, $false, $false, $true, 1, $false, "12×59=", 2) | Out-Null
$d.Content.Find.Execute("92×85=", $true, $false, $false, $false, $false, $true, 1, $false, "49×53=", 2) | Out-Null
$d.Content.Find.Execute("45×36=", $true, $false, $false, $false, $false, $true, 1, $false, "17×100=", 2) | Out-Null
$d.Content.Find.Execute("79×60=", $true, $false, $false, $false, $false, $true, 1, $false, "42×77=", 2) | Out-Null
$d.Content.Find.Execute("35×19=", $true, $false, $false, $false, $false, $true, 1, $false, "20×62=", 2) | Out-Null
$d.Content.Find.Execute("84×26=", $true, $false, $false, $false, $false, $true, 1, $false, "42×55=", 2) | Out-Null
$d.Content.Find.Execute("10×34=", $true, $false, $false, $false, $false, $true, 1, $false, "25×62=", 2) | Out-Null
$d.Content.Find.Execute("12×16=", $true, $false, $false, $false, $false, $true, 1, $false, "86×34=", 2) | Out-Null
$d.Content.Find.Execute("10×95=", $true, $false, $false, $false, $false, $true, 1, $false, "18×10=", 2) | Out-Null
$d.Content.Find.Execute("53×79=", $true, $false, $false, $false, $false, $true, 1, $false, "44×80=", 2) | Out-Null
$d.Content.Find.Execute("81×24=", $true, $false, $false, $false, $false, $true, 1, $false, "26×50=", 2) | Out-Null
$d.Content.Find.Execute("44×43=", $true, $false, $false, $false, $false, $true, 1, $false, "21×79=", 2) | Out-Null
$d.Content.Find.Execute("45×66=", $true, $false, $false, $false, $false, $true, 1, $false, "34×40=", 2) | Out-Null
$d.Content.Find.Execute("36×93=", $true, $false, $false, $false, $false, $true, 1, $false, "97×30=", 2) | Out-Null
$d.Content.Find.Execute("50×28=", $true, $false, $false, $false, $false, $true, 1, $false, "29×36=", 2) | Out-Null
$d.Content.Find.Execute("96×98=", $true, $false, $false, $false, $false, $true, 1, $false, "25×100=", 2) | Out-Null
$d.Content.Find.Execute("41×89=", $true, $false, $false, $false, $false, $true, 1, $false, "45×12=", 2) | Out-Null
$d.Content.Find.Execute("79×96=", $true, $false, $false, $false, $false, $true, 1, $false, "47×41=", 2) | Out-Null
$d.Content.Find.Execute("10×100=", $true, $false, $false, $false, $false, $true, 1, $false, "95×46=", 2) | Out-Null
$d.Content.Find.Execute("63×69=", $true, $false, $false, $false, $false, $true, 1, $false, "65×82=", 2) | Out-Null
$d.Content.Find.Execute("54×65=", $true, $false, $false, $false, $false, $true, 1, $false, "59×30=", 2) | Out-Null
